$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

# Translate frequency-related categorical values from French to English
$used.Replace("rarement", "Rarely")
$used.Replace("occasionnellement", "Occasionally")
$used.Replace("Fréquente", "Frequently")

# Translate chronic-disease-related categorical values from French to English
$used.Replace("autres", "Others")
$used.Replace("aucune maladie", "None")
$used.Replace("diabete", "Diabetes")
$used.Replace("HTA", "Hypertension")
$used.Replace("asthme", "Asthma")
